# Refresh '想去人数' (interest count, column F) figures to the latest scrape snapshot
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1854
$ws.Cells.Item(3, 6).Value = 21
$ws.Cells.Item(5, 6).Value = 54
$ws.Cells.Item(8, 6).Value = 199
$ws.Cells.Item(9, 6).Value = 643
$ws.Cells.Item(10, 6).Value = 69
$ws.Cells.Item(12, 6).Value = 765
$ws.Cells.Item(13, 6).Value = 1461
$ws.Cells.Item(14, 6).Value = 1218
$ws.Cells.Item(15, 6).Value = 1447
$ws.Cells.Item(16, 6).Value = 26
$ws.Cells.Item(17, 6).Value = 1255
$ws.Cells.Item(19, 6).Value = 1592
$ws.Cells.Item(21, 6).Value = 1017
$ws.Cells.Item(22, 6).Value = 326
$ws.Cells.Item(25, 6).Value = 1400
$ws.Cells.Item(26, 6).Value = 93
$ws.Cells.Item(29, 6).Value = 1075
$ws.Cells.Item(30, 6).Value = 263442
$ws.Cells.Item(31, 6).Value = 978
$ws.Cells.Item(32, 6).Value = 19
$ws.Cells.Item(33, 6).Value = 556
$ws.Cells.Item(34, 6).Value = 1322
$ws.Cells.Item(35, 6).Value = 1044
$ws.Cells.Item(36, 6).Value = 887
$ws.Cells.Item(37, 6).Value = 1071
$ws.Cells.Item(38, 6).Value = 22
$ws.Cells.Item(39, 6).Value = 51
$ws.Cells.Item(40, 6).Value = 42
$ws.Cells.Item(41, 6).Value = 837
$ws.Cells.Item(42, 6).Value = 1598
$ws.Cells.Item(43, 6).Value = 93
$ws.Cells.Item(44, 6).Value = 39
$ws.Cells.Item(45, 6).Value = 794
$ws.Cells.Item(46, 6).Value = 87
$ws.Cells.Item(47, 6).Value = 777

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 41
$ws.Cells.Item(3, 6).Value = 49
$ws.Cells.Item(4, 6).Value = 133
$ws.Cells.Item(10, 6).Value = 174
$ws.Cells.Item(11, 6).Value = 1437
$ws.Cells.Item(12, 6).Value = 72
$ws.Cells.Item(14, 6).Value = 2535
$ws.Cells.Item(15, 6).Value = 1190
$ws.Cells.Item(16, 6).Value = 395
$ws.Cells.Item(18, 6).Value = 223
$ws.Cells.Item(20, 6).Value = 68
$ws.Cells.Item(26, 6).Value = 285
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 15
$ws.Cells.Item(32, 6).Value = 234
$ws.Cells.Item(34, 6).Value = 140
$ws.Cells.Item(39, 6).Value = 167
$ws.Cells.Item(42, 6).Value = 16
$ws.Cells.Item(43, 6).Value = 41
$ws.Cells.Item(44, 6).Value = 41
$ws.Cells.Item(46, 6).Value = 125
$ws.Cells.Item(47, 6).Value = 55

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 248
$ws.Cells.Item(5, 6).Value = 2778
$ws.Cells.Item(6, 6).Value = 4528
$ws.Cells.Item(8, 6).Value = 6
$ws.Cells.Item(9, 6).Value = 536
$ws.Cells.Item(10, 6).Value = 651
$ws.Cells.Item(11, 6).Value = 429
$ws.Cells.Item(12, 6).Value = 224
$ws.Cells.Item(13, 6).Value = 813
$ws.Cells.Item(14, 6).Value = 203
$ws.Cells.Item(15, 6).Value = 465

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1854
$ws.Cells.Item(3, 6).Value = 248
$ws.Cells.Item(4, 6).Value = 2778
$ws.Cells.Item(5, 6).Value = 4528
$ws.Cells.Item(6, 6).Value = 651
$ws.Cells.Item(7, 6).Value = 54
$ws.Cells.Item(8, 6).Value = 224
$ws.Cells.Item(9, 6).Value = 224
$ws.Cells.Item(10, 6).Value = 813
$ws.Cells.Item(11, 6).Value = 813
$ws.Cells.Item(12, 6).Value = 203
$ws.Cells.Item(15, 6).Value = 199
$ws.Cells.Item(16, 6).Value = 1437
$ws.Cells.Item(18, 6).Value = 765
$ws.Cells.Item(19, 6).Value = 2535
$ws.Cells.Item(20, 6).Value = 1190
$ws.Cells.Item(21, 6).Value = 1461
$ws.Cells.Item(22, 6).Value = 1218
$ws.Cells.Item(23, 6).Value = 1447
$ws.Cells.Item(24, 6).Value = 1255
$ws.Cells.Item(25, 6).Value = 223
$ws.Cells.Item(26, 6).Value = 68
$ws.Cells.Item(27, 6).Value = 1592
$ws.Cells.Item(29, 6).Value = 1017
$ws.Cells.Item(30, 6).Value = 326
$ws.Cells.Item(31, 6).Value = 465
$ws.Cells.Item(32, 6).Value = 465
$ws.Cells.Item(34, 6).Value = 1400
$ws.Cells.Item(37, 6).Value = 1075
$ws.Cells.Item(38, 6).Value = 285
$ws.Cells.Item(39, 6).Value = 978
$ws.Cells.Item(40, 6).Value = 19
$ws.Cells.Item(41, 6).Value = 1044
$ws.Cells.Item(42, 6).Value = 887
$ws.Cells.Item(43, 6).Value = 1071
$ws.Cells.Item(45, 6).Value = 837
$ws.Cells.Item(47, 6).Value = 1598
$ws.Cells.Item(48, 6).Value = 93
$ws.Cells.Item(49, 6).Value = 794
$ws.Cells.Item(50, 6).Value = 41
$ws.Cells.Item(51, 6).Value = 777
$ws.Cells.Item(53, 6).Value = 55
